$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet per the diff (map.contact_board -> map.contact_education)
$ws.Name = "map.contact_education"

# Remove now-obsolete transformation-check rows (9-11); table shrinks to A1:M8
$ws.Rows("9:11").Delete()

# Row 2
$v1 = @'
TC01_Onbaording the map.contact_education   into MDM Raw layer.
'@
$ws.Range("C2").Value = $v1
$v2 = @'
TC_Verification of structure for onbaordingmap.contact_education
'@
$ws.Range("D2").Value = $v2
$v3 = @'
 Open the Microsoft SQL Server,Go to the Databases folder and expand it. 
                Go to the Tables folder right click on the map.contact_educationtable a new window will open with metadata. 
                Go to the Tables folder right click on the map.contact_educationtable a new window will open with metadata. 
                Verify the struture of the both the tables. 
'@
$ws.Range("F2").Value = $v3
$v4 = @'
Structure should match between Source and Target table.
'@
$ws.Range("G2").Value = $v4
$v5 = @'
Structure should match between Source and Target table.
'@
$ws.Range("H2").Value = $v5

# Row 3
$v6 = @'
TC00_Countvalidation_map.contact_education
'@
$ws.Range("C3").Value = $v6
$v7 = @'
TC_Verification of count between the source and target table
'@
$ws.Range("D3").Value = $v7
$v8 = @'
 Verify the count of the source table using below query in the MDM Raw Layer 
 select count(*) from crm.s_contact where change_flag= 'Y' 
                   Verify the count of the target table using below query in the MDM MAP layer 
 select count(*) from map.contact_education 
                   Note: Please apply processID filter for both the tables Verify the count between the  source and target table. 
'@
$ws.Range("F3").Value = $v8
$v9 = @'
Count should match between source and target on the latest process ID.
'@
$ws.Range("G3").Value = $v9
$v10 = @'
Count should match between source and target on the latest process ID.
'@
$ws.Range("H3").Value = $v10

# Row 4
$v11 = @'
TC_Verification of duplicates in the data
'@
$ws.Range("C4").Value = $v11
$v12 = @'
Check whether there  any duplicates exists on the latest processID in map.contact_educationtable.
'@
$ws.Range("D4").Value = $v12
$v13 = @'
Check whether there are any duplicates in the data loaded on the latest processID
'@
$ws.Range("F4").Value = $v13
$v14 = @'
 Below Query will be used to identify the duplicates  
                select mdmid,count(*) from map.contact_educationwhere 1=1 group by mdmid having count(*) 
                Note: Please apply processID filter while validating the data  
'@
$ws.Range("G4").Value = $v14
$v15 = @'
 Duplicates should not exists on the latest processID. We should have all unique records 
'@
$ws.Range("H4").Value = $v15

# Row 5
$v16 = @'
TC_TransformationCheck_Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("C5").Value = $v16
$v17 = @'
TC_Verify whether the transformation logic has applied as per the requirement document for the column  Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("D5").Value = $v17
$v18 = @'
Extract the data from the table in the MDM Raw Layer using  below query
SourceQuery: 'select row_id,'Degree' as EducationTypeCode from crm.s_Contact where 1=1 and c1.ROW_ID IN (SELECT ROW_ID FROM CRM.CONTACTID_STG WHERE PROCESSID='$LatestProcessID')  and change_flag = 'Y' Extract the data from the table in the MDM Map Layer using  below query 
 TargetQuery: select row_id,EducationTypeCode   from map.contact_education Data Validation:  Verify the data from Source Query and Target Query  i.e by applying the except query.   
 Note: Please apply processID filter while validating the data for both the above mentioned tables.
'@
$ws.Range("F5").Value = $v18
$v19 = @'
Data should be loaded as per the Transformation logic
'@
$ws.Range("G5").Value = $v19
$v20 = @'
Extract the data from the table in the MDM Map Layer using  below query  value should match between the Source and Target tables
'@
$ws.Range("H5").Value = $v20

# Row 6
$v21 = @'
TC_TransformationCheck_Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("C6").Value = $v21
$v22 = @'
TC_Verify whether the transformation logic has applied as per the requirement document for the column  Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("D6").Value = $v22
$v23 = @'
Extract the data from the table in the MDM Raw Layer using  below query
SourceQuery: select row_id, CASE WHEN c1.DEGREE IS NULL THEN 'Not Provided' ELSE SRC1.DEGREE END EducationDegree from crm.s_Contact where 1=1 and c1.ROW_ID IN (SELECT ROW_ID FROM CRM.CONTACTID_STG WHERE PROCESSID='$LatestProcessID')  and change_flag = 'Y' Extract the data from the table in the MDM Map Layer using  below query 
 TargetQuery: select row_id,EducationDegree  from map.contact_education Data Validation:  Verify the data from Source Query and Target Query  i.e by applying the except query.   
 Note: Please apply processID filter while validating the data for both the above mentioned tables.
'@
$ws.Range("F6").Value = $v23
$v24 = @'
Data should be loaded as per the Transformation logic
'@
$ws.Range("G6").Value = $v24
$v25 = @'
Extract the data from the table in the MDM Map Layer using  below query  value should match between the Source and Target tables
'@
$ws.Range("H6").Value = $v25

# Row 7
$v26 = @'
TC_TransformationCheck_Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("C7").Value = $v26
$v27 = @'
TC_Verify whether the transformation logic has applied as per the requirement document for the column  Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("D7").Value = $v27
$v28 = @'
Extract the data from the table in the MDM Raw Layer using  below query
SourceQuery: select row_id, c2.ATTRIB_04 as EducationBackground  from crm.s_Contact c1, crm.s_contact_x c2 where 1=1 and c1.row_id = c2.row_id  and c1.ROW_ID IN (SELECT ROW_ID FROM CRM.CONTACTID_STG WHERE PROCESSID='$LatestProcessID')  and change_flag = 'Y' Extract the data from the table in the MDM Map Layer using  below query 
 TargetQuery: select row_id,EducationBackground   from map.contact_education Data Validation:  Verify the data from Source Query and Target Query  i.e by applying the except query.   
 Note: Please apply processID filter while validating the data for both the above mentioned tables.
'@
$ws.Range("F7").Value = $v28
$v29 = @'
Data should be loaded as per the Transformation logic
'@
$ws.Range("G7").Value = $v29
$v30 = @'
Extract the data from the table in the MDM Map Layer using  below query  value should match between the Source and Target tables
'@
$ws.Range("H7").Value = $v30

# Row 8
$v31 = @'
TC_TransformationCheck_Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("C8").Value = $v31
$v32 = @'
TC_Verify whether the transformation logic has applied as per the requirement document for the column  Extract the data from the table in the MDM Map Layer using  below query
'@
$ws.Range("D8").Value = $v32
$v33 = @'
Extract the data from the table in the MDM Raw Layer using  below query
SourceQuery: select row_id, c2.X_ALUM_HON_QUAL AS Qualifications  from crm.s_Contact c1, crm.s_contact_x c2 where 1=1 and c1.row_id = c2.row_id  and c1.ROW_ID IN (SELECT ROW_ID FROM CRM.CONTACTID_STG WHERE PROCESSID='$LatestProcessID')  and change_flag = 'Y' Extract the data from the table in the MDM Map Layer using  below query 
 TargetQuery: select row_id,Qualifications    from map.contact_education Data Validation:  Verify the data from Source Query and Target Query  i.e by applying the except query.   
 Note: Please apply processID filter while validating the data for both the above mentioned tables.
'@
$ws.Range("F8").Value = $v33
$v34 = @'
Data should be loaded as per the Transformation logic
'@
$ws.Range("G8").Value = $v34
$v35 = @'
Extract the data from the table in the MDM Map Layer using  below query  value should match between the Source and Target tables
'@
$ws.Range("H8").Value = $v35
